$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column N (14) to column O (15) for rows 3-14,
# then set the new values in column O.

$ws.Range("N3:N14").Copy() | Out-Null
$ws.Range("O3:O14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Row 3: empty cell (format only, already copied above)

# Row 4: header year
$ws.Cells.Item(4, 15).Value = 2021

# Data rows 5-14
$ws.Cells.Item(5, 15).Value = 70.636215334420882
$ws.Cells.Item(6, 15).Value = 107.1
$ws.Cells.Item(7, 15).Value = 55.452054794520542
$ws.Cells.Item(8, 15).Value = 84.375
$ws.Cells.Item(9, 15).Value = 120.48192771084337
$ws.Cells.Item(10, 15).Value = 109.53346855983774
$ws.Cells.Item(11, 15).Value = 147.7690288713911
$ws.Cells.Item(12, 15).Value = 25.545675020210183
$ws.Cells.Item(13, 15).Value = 82.457854874175425
$ws.Cells.Item(14, 15).Value = 15.384615384615385
